# Finish LPs and OMs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordered list of LP values that should occupy column A (rows 1-21)
$values = @(
    "LP",
    "LP-032383",
    "LP-043257",
    "LP-044468",
    "LP-047127",
    "LP-047276",
    "LP-047425",
    "LP-047636",
    "LP-048244",
    "LP-048301",
    "LP-048670",
    "LP-048674",
    "LP-048746",
    "LP-048932",
    "LP-049257",
    "LP-049427",
    "LP-049443",
    "LP-049800",
    "LP-050182",
    "LP-050183",
    "LP-050208"
)

# Determine how many rows currently contain data in column A so we know
# whether we need to clear any trailing rows (old data had 24 rows, new
# data only needs 21).
$usedRange = $ws.UsedRange
$oldLastRow = $usedRange.Rows.Count

# Write the new values into A1:A21
for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

# Clear any now-unused trailing rows (rows 22-24 previously held data)
$newLastRow = $values.Count
if ($oldLastRow -gt $newLastRow) {
    $clearRange = $ws.Range($ws.Cells.Item($newLastRow + 1, 1), $ws.Cells.Item($oldLastRow, 1))
    $clearRange.ClearContents()
}

# Remove the previous selection on the sheet view (selection A2:A24 -> none specific)
$ws.Range("A1").Select()
